$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Split the A2:A5 "Minnow3 Module" merge into two 2-row groups -------
$ws.Range("A2:A5").UnMerge()

# A3 should pick up the "bottom of group" border (same as A5 currently has)
$ws.Range("A5").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# A4 should pick up the "top of group" border (same as A2 currently has)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. New label text for the two remaining build groups -------------------
$ws.Range("A2").Value = "Leaf Hill"
$ws.Range("B2").Value = "FAB D"
$ws.Range("A4").Value = "UP2"
$ws.Range("B4").Value = "FAB A"

# --- 3. Re-merge the label columns into the new 2-row groups ---------------
$ws.Range("A2:A3").Merge()
$ws.Range("A4:A5").Merge()

# --- 4. Add the new "Disable flash region access(R)" column (H) ------------
$ws.Columns("H").ColumnWidth = 27

$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Disable flash region access(R)"

$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "N/A"

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = "N/A"

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = "N/A"

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "Y"

# --- 5. Drop the "Leaf Hill" / "UP2" rows that used to live further down ---
$ws.Rows("6:9").Delete()

# --- 6. Restore the selection -----------------------------------------------
$ws.Range("D5").Select()
